# Update "想去人数" (number of interested attendees) figures that were
# refreshed when the site was regenerated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): rows 2-10 map directly to the updated counts.
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 317
$ws1.Range("F3").Value = 63
$ws1.Range("F4").Value = 479
$ws1.Range("F5").Value = 4636
$ws1.Range("F6").Value = 362
$ws1.Range("F7").Value = 626
$ws1.Range("F9").Value = 718
$ws1.Range("F10").Value = 200

# Sheet "全部类型" (sheet4): same events, but shifted down one row from
# row 10 onward because this sheet contains an extra event row.
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 317
$ws4.Range("F3").Value = 63
$ws4.Range("F4").Value = 479
$ws4.Range("F5").Value = 4636
$ws4.Range("F6").Value = 362
$ws4.Range("F7").Value = 626
$ws4.Range("F9").Value = 718
$ws4.Range("F11").Value = 200
